$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row above the old "Sum" row (row 16). Excel shifts the Sum
# row (with its formulas/formatting) down to row 17, and the new row 16
# inherits formatting from the row above it.
# ---------------------------------------------------------------------------
$ws.Rows("16:16").Insert() | Out-Null

# Row 15 previously had no Status (column G) value - it now gets one.
$ws.Range("G14").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").Value = "zu bestellen"

# ---------------------------------------------------------------------------
# New data row 16: "silicon cable" purchase entry.
# Add the hyperlink first so the URL's shared-string slot is created before
# the "silicon cable" label's slot (matches authoring order).
# ---------------------------------------------------------------------------
$url = "https://www.amazon.de/-/en/LAVMHAB-Silicone-Flexible-Pumping-Transfer/dp/B0D9NRLR1K?crid=2O9NW9O1OLT95&dib=eyJ2IjoiMSJ9.FCrh94MEEO0jp4-k4fQGxYMMit5zN8Ddzby2BDyEk_UCjtk5-yg4-0wtAVgOpdvtJPhVu5nHl7yQQ_OHwYc5ZDKkRoLd3MJQr3SzEbZ9XvGr4whdEfGauaxCsqcwoGINfLkiOMZ8UWXYtKcGZ75nICn1K4kxyaDTrK6Kbp4VFcNp2QmsZN7TEfri5icLJBDTrbwDlcSHkvuQJmqD_u3aviGNadizMwUQCy1ILMNom7oc3HvyXUFJgnUfMLOP75BCp-aI-RKrtz4vaYyVyocRhJOfuYPtpJOi2f2kxZu0wzo.3ozZ4Wradsj81AlhWYepdDAu5nTMiNkJFeO6_YZtgps&dib_tag=se&keywords=silikonschlauch%2B5mm&nsdOptOutParam=true&qid=1732518622&sprefix=5mm%2Bsilic%2Caps%2C109&sr=8-17&th=1"

$ws.Hyperlinks.Add($ws.Range("B16"), $url, $null, $null, $url) | Out-Null
# Re-apply the same "hyperlink cell" formatting used by the other link cells
# (B9 ... ) so the style index matches the rest of the table.
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null

# A16: product label, formatted like the other product-name cells (A14 ...).
$ws.Range("A14").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = "silicon cable"

# Price + split columns.
$ws.Range("C16").Value = 10.07
$ws.Range("D16").Formula = "=C16/3"
$ws.Range("E16").Formula = "=C16/3"
$ws.Range("F16").Formula = "=C16/3"

# Status column.
$ws.Range("G16").Value = "zu bestellen"

# ---------------------------------------------------------------------------
# Row 17 (formerly row 16) is the "Sum" row - extend the price sum to cover
# the newly inserted row.
# ---------------------------------------------------------------------------
$ws.Range("C17").Formula = "=SUM(C7:C16)"

# ---------------------------------------------------------------------------
# View: zoom + selected cell.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 200 | Out-Null
$ws.Range("E21").Select() | Out-Null
